$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Value = "[60.2473512910275, 67.31262613449297]"
$ws.Range("T2").Value = "[47.653811673308354, 52.35145618165127]"
$ws.Range("L3").Value = "[59.33879702384884, 67.8875825043677]"
$ws.Range("T3").Value = "[47.807431269157206, 52.216712100026086]"
